$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header in column J
$ws.Range("J1").Value = "path_version"

# Fill J2:J29 with the path_version value (2) for each existing data row
for ($r = 2; $r -le 29; $r++) {
    $ws.Cells.Item($r, 10).Value = 2
}

# Move the active selection to K27, matching the post-edit cursor position
$ws.Range("K27").Select()
